# Add two new controlled-vocabulary instrument model values to the
# ENA experiment metadata template:
#   - "DNBSEQ-T10x4RS"          (inserted after "DNBSEQ-G50")
#   - "Illumina NovaSeq X Plus" (inserted after "Illumina NovaSeq X")
#
# The instrument-model pick list lives on the hidden 'cv_experiment'
# sheet, column M, and is referenced by the 'instrumentmodel' workbook
# defined name (used for data validation on the 'experiment' sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cv_experiment")

# Full, alphabetically-sorted instrument-model list, with the two new
# entries already inserted in their correct sorted position.
$models = @(
    '454 GS',
    '454 GS 20',
    '454 GS FLX',
    '454 GS FLX Titanium',
    '454 GS FLX+',
    '454 GS Junior',
    'AB 310 Genetic Analyzer',
    'AB 3130 Genetic Analyzer',
    'AB 3130xL Genetic Analyzer',
    'AB 3500 Genetic Analyzer',
    'AB 3500xL Genetic Analyzer',
    'AB 3730 Genetic Analyzer',
    'AB 3730xL Genetic Analyzer',
    'AB 5500 Genetic Analyzer',
    'AB 5500xl Genetic Analyzer',
    'AB 5500xl-W Genetic Analysis System',
    'AB SOLiD 3 Plus System',
    'AB SOLiD 4 System',
    'AB SOLiD 4hq System',
    'AB SOLiD PI System',
    'AB SOLiD System',
    'AB SOLiD System 2.0',
    'AB SOLiD System 3.0',
    'BGISEQ-50',
    'BGISEQ-500',
    'Complete Genomics',
    'DNBSEQ-G400',
    'DNBSEQ-G400 FAST',
    'DNBSEQ-G50',
    'DNBSEQ-T10x4RS',
    'DNBSEQ-T7',
    'Element AVITI',
    'FASTASeq 300',
    'GENIUS',
    'GS111',
    'Genapsys Sequencer',
    'GenoCare 1600',
    'GenoLab M',
    'GridION',
    'Helicos HeliScope',
    'HiSeq X Five',
    'HiSeq X Ten',
    'Illumina Genome Analyzer',
    'Illumina Genome Analyzer II',
    'Illumina Genome Analyzer IIx',
    'Illumina HiScanSQ',
    'Illumina HiSeq 1000',
    'Illumina HiSeq 1500',
    'Illumina HiSeq 2000',
    'Illumina HiSeq 2500',
    'Illumina HiSeq 3000',
    'Illumina HiSeq 4000',
    'Illumina HiSeq X',
    'Illumina MiSeq',
    'Illumina MiniSeq',
    'Illumina NovaSeq 6000',
    'Illumina NovaSeq X',
    'Illumina NovaSeq X Plus',
    'Illumina iSeq 100',
    'Ion GeneStudio S5',
    'Ion GeneStudio S5 Plus',
    'Ion GeneStudio S5 Prime',
    'Ion Torrent Genexus',
    'Ion Torrent PGM',
    'Ion Torrent Proton',
    'Ion Torrent S5',
    'Ion Torrent S5 XL',
    'MGISEQ-2000RS',
    'MinION',
    'NextSeq 1000',
    'NextSeq 2000',
    'NextSeq 500',
    'NextSeq 550',
    'Onso',
    'PacBio RS',
    'PacBio RS II',
    'PromethION',
    'Revio',
    'Sentosa SQ301',
    'Sequel',
    'Sequel II',
    'Sequel IIe',
    'Tapestri',
    'UG 100',
    'unspecified'
)

$count = $models.Count

# Write the whole column in one shot as a vertical array.
$arr = New-Object 'object[,]' $count,1
for ($i = 0; $i -lt $count; $i++) {
    $arr[$i, 0] = $models[$i]
}
$ws.Range($ws.Cells.Item(1, 13), $ws.Cells.Item($count, 13)).Value2 = $arr

# Update the 'instrumentmodel' defined name so the data-validation list
# on the 'experiment' sheet covers the new, longer range.
$name = $wb.Names.Item("instrumentmodel")
$name.RefersTo = "=cv_experiment!`$M`$1:`$M`$" + $count
